$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N2").Value = -475.5
$ws.Range("M2").Value = -61.875
$ws.Range("L2").Value = 249.5
$ws.Range("J2").Value = 249.5
$ws.Range("H2").Value = 189.8
$ws.Range("K2").Value = 174.875
$ws.Range("I2").Value = 174.875
$ws.Range("H15").Value = 853.0769
$ws.Range("M15").Value = -2390.2307
$ws.Range("K15").Value = 2559.2307
$ws.Range("I15").Value = 853.0769
$ws.Range("L64").Value = 3949
$ws.Range("N64").Value = -4445
$ws.Range("J64").Value = 3949
$ws.Range("H64").Value = 4870.5
$ws.Range("L67").Value = 3949
$ws.Range("N67").Value = -5665
$ws.Range("J67").Value = 3949
$ws.Range("H67").Value = 4870.5
$ws.Range("L86").Value = 2485.6
$ws.Range("N86").Value = -4731.6
$ws.Range("M86").Value = -2334.0715
$ws.Range("H86").Value = 3201.4211
$ws.Range("J86").Value = 2485.6
$ws.Range("K86").Value = 3457.0715
$ws.Range("I86").Value = 3457.0715
$ws.Range("M89").Value = -11669.3575
$ws.Range("L89").Value = 12428
$ws.Range("N89").Value = -23660
$ws.Range("J89").Value = 2485.6
$ws.Range("H89").Value = 3201.4211
$ws.Range("K89").Value = 17285.3575
$ws.Range("I89").Value = 3457.0715
$ws.Range("H100").Value = 2468.7144
$ws.Range("M100").Value = -2015.077
$ws.Range("K100").Value = 2556.077
$ws.Range("I100").Value = 2556.077
$ws.Range("K101").Value = 2246.625
$ws.Range("I101").Value = 748.875
$ws.Range("N101").Value = -4608.25
$ws.Range("M101").Value = -624.625
$ws.Range("L101").Value = 1364.25
$ws.Range("J101").Value = 454.75
$ws.Range("H101").Value = 650.8333
$ws.Range("M137").Value = -76923816
$ws.Range("H137").Value = 16668105
$ws.Range("K137").Value = 76926366
$ws.Range("I137").Value = 25642122
$ws.Range("M141").Value = -3870.25
$ws.Range("H141").Value = 2817.3
$ws.Range("K141").Value = 9050.25
$ws.Range("I141").Value = 3016.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K11").Value = 8573214
$ws.Range("I11").Value = 8573214
$ws.Range("M11").Value = -8573070
$ws.Range("H11").Value = 8573214
$ws.Range("M61").Value = -3365.2856
$ws.Range("H61").Value = 3577.2856
$ws.Range("K61").Value = 3577.2856
$ws.Range("I61").Value = 3577.2856
$ws.Range("L74").Value = 3372.4
$ws.Range("N74").Value = -5120.4
$ws.Range("J74").Value = 3372.4
$ws.Range("H74").Value = 2107.9412
$ws.Range("M74").Value = -707.0834
$ws.Range("K74").Value = 1581.0834
$ws.Range("I74").Value = 1581.0834
$ws.Range("M77").Value = -3537.416999999999
$ws.Range("L77").Value = 16862
$ws.Range("N77").Value = -25598
$ws.Range("J77").Value = 3372.4
$ws.Range("H77").Value = 2107.9412
$ws.Range("K77").Value = 7905.416999999999
$ws.Range("I77").Value = 1581.0834
$ws.Range("K107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 67962.336
$ws.Range("L107").Value = 67962.336
$ws.Range("N107").Value = -75642.336
$ws.Range("H107").Value = 67962.336
$ws.Range("K136").Value = 10731.8568
$ws.Range("I136").Value = 3577.2856
$ws.Range("M136").Value = -8181.856800000001
$ws.Range("H136").Value = 3577.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M16").Value = -446.3333
$ws.Range("L16").Value = 754.5
$ws.Range("J16").Value = 754.5
$ws.Range("H16").Value = 671.6
$ws.Range("K16").Value = 616.3333
$ws.Range("I16").Value = 616.3333
$ws.Range("N16").Value = -1094.5
$ws.Range("L22").Value = 748.25
$ws.Range("N22").Value = -1094.25
$ws.Range("M22").Value = -360.2222
$ws.Range("J22").Value = 748.25
$ws.Range("H22").Value = 599.38464
$ws.Range("K22").Value = 533.2222
$ws.Range("I22").Value = 533.2222
$ws.Range("J55").Value = 102166.336
$ws.Range("H55").Value = 102166.336
$ws.Range("L55").Value = 102166.336
$ws.Range("N55").Value = -102712.336
$ws.Range("L86").Value = 30330248
$ws.Range("N86").Value = -30332494
$ws.Range("M86").Value = -20403.066
$ws.Range("H86").Value = 12844447
$ws.Range("J86").Value = 30330248
$ws.Range("K86").Value = 21526.066
$ws.Range("I86").Value = 21526.066
$ws.Range("M89").Value = -102014.33
$ws.Range("L89").Value = 151651240
$ws.Range("N89").Value = -151662472
$ws.Range("J89").Value = 30330248
$ws.Range("H89").Value = 12844447
$ws.Range("K89").Value = 107630.33
$ws.Range("I89").Value = 21526.066
$ws.Range("K107").Value = 809
$ws.Range("I107").Value = 809
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1111
$ws.Range("H107").Value = 809

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J4").Value = 0
$ws.Range("H4").Value = 16345
$ws.Range("L4").Value = 0
$ws.Range("L7").Value = 37.5
$ws.Range("H7").Value = 50000176
$ws.Range("K7").Value = 55555748
$ws.Range("I7").Value = 55555748
$ws.Range("N7").Value = -263.5
$ws.Range("J7").Value = 37.5
$ws.Range("M7").Value = -55555635
$ws.Range("M16").Value = 27.57144
$ws.Range("L16").Value = 144
$ws.Range("J16").Value = 144
$ws.Range("H16").Value = 245
$ws.Range("K16").Value = 259.42856
$ws.Range("I16").Value = 259.42856
$ws.Range("N16").Value = -718
$ws.Range("M58").Value = -2259.8096
$ws.Range("H58").Value = 2233.9666
$ws.Range("K58").Value = 2462.8096
$ws.Range("I58").Value = 2462.8096
$ws.Range("K62").Value = 2412.25
$ws.Range("I62").Value = 2412.25
$ws.Range("M62").Value = -1788.25
$ws.Range("H62").Value = 2441.3333
$ws.Range("H65").Value = 2441.3333
$ws.Range("K65").Value = 12061.25
$ws.Range("I65").Value = 2412.25
$ws.Range("M65").Value = -8941.25
$ws.Range("N92").Value = -38992
$ws.Range("L92").Value = 34000
$ws.Range("J92").Value = 34000
$ws.Range("H92").Value = 34000
$ws.Range("H105").Value = 985.9
$ws.Range("M105").Value = 707.1111000000001
$ws.Range("K105").Value = 1039.8889
$ws.Range("I105").Value = 1039.8889
$ws.Range("K107").Value = 266.77777
$ws.Range("I107").Value = 266.77777
$ws.Range("J107").Value = 250
$ws.Range("L107").Value = 250
$ws.Range("N107").Value = -4090
$ws.Range("M107").Value = 1653.22223
$ws.Range("H107").Value = 263.72726
$ws.Range("J113").Value = 144
$ws.Range("H113").Value = 245
$ws.Range("K113").Value = 259.42856
$ws.Range("I113").Value = 259.42856
$ws.Range("N113").Value = -4484
$ws.Range("M113").Value = 1910.57144
$ws.Range("L113").Value = 144
$ws.Range("L122").Value = 8995.5
$ws.Range("N122").Value = -13895.5
$ws.Range("J122").Value = 2998.5
$ws.Range("H122").Value = 2997.3333
$ws.Range("M134").Value = -6545.571599999999
$ws.Range("L134").Value = 5997
$ws.Range("N134").Value = -11067
$ws.Range("J134").Value = 1999
$ws.Range("H134").Value = 2898.375
$ws.Range("K134").Value = 9080.5716
$ws.Range("I134").Value = 3026.8572
$ws.Range("K136").Value = 7388.4288
$ws.Range("I136").Value = 2462.8096
$ws.Range("M136").Value = -4838.4288
$ws.Range("H136").Value = 2233.9666
$ws.Range("L141").Value = 397034.9
$ws.Range("N141").Value = -407394.9
$ws.Range("J141").Value = 397034.9
$ws.Range("H141").Value = 397034.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N2").Value = -4800371.199999999
$ws.Range("L2").Value = 4800145.199999999
$ws.Range("J2").Value = 800024.2
$ws.Range("H2").Value = 588268.5
$ws.Range("H9").Value = 4833.3335
$ws.Range("M16").Value = -2528.5
$ws.Range("H16").Value = 900.5
$ws.Range("K16").Value = 2701.5
$ws.Range("I16").Value = 900.5
$ws.Range("L22").Value = 11250
$ws.Range("N22").Value = -11588
$ws.Range("J22").Value = 3750
$ws.Range("H22").Value = 3750
$ws.Range("L27").Value = 11250
$ws.Range("N27").Value = -11454
$ws.Range("J27").Value = 3750
$ws.Range("H27").Value = 3750
$ws.Range("H33").Value = 158.90909
$ws.Range("J33").Value = 193.33333
$ws.Range("K33").Value = 876
$ws.Range("I33").Value = 146
$ws.Range("N33").Value = -1725.99998
$ws.Range("M33").Value = -593
$ws.Range("L33").Value = 1159.99998
$ws.Range("L41").Value = 9006
$ws.Range("N41").Value = -9682
$ws.Range("J41").Value = 3002
$ws.Range("H41").Value = 3002
$ws.Range("J68").Value = 796
$ws.Range("H68").Value = 919
$ws.Range("K68").Value = 2941.5
$ws.Range("I68").Value = 980.5
$ws.Range("N68").Value = -4010
$ws.Range("M68").Value = -2130.5
$ws.Range("L68").Value = 2388
$ws.Range("J71").Value = 796
$ws.Range("L71").Value = 7164
$ws.Range("M71").Value = -4768.5
$ws.Range("H71").Value = 919
$ws.Range("K71").Value = 8824.5
$ws.Range("I71").Value = 980.5
$ws.Range("N71").Value = -15276
$ws.Range("M132").Value = -9620
$ws.Range("H132").Value = 1366.6666
$ws.Range("K132").Value = 12150
$ws.Range("I132").Value = 1350
$ws.Range("M138").Value = -31866.33199999999
$ws.Range("H138").Value = 12191.9
$ws.Range("K138").Value = 37006.33199999999
$ws.Range("I138").Value = 12335.444

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M70").Value = -3394
$ws.Range("H70").Value = 3998
$ws.Range("K70").Value = 3664
$ws.Range("I70").Value = 3664
$ws.Range("M73").Value = -2728
$ws.Range("H73").Value = 3998
$ws.Range("K73").Value = 3664
$ws.Range("I73").Value = 3664
$ws.Range("N92").Value = -17354.5
$ws.Range("L92").Value = 13610.5
$ws.Range("J92").Value = 13610.5
$ws.Range("H92").Value = 13610.5
$ws.Range("M132").Value = -5153.428400000001
$ws.Range("H132").Value = 15153967
$ws.Range("K132").Value = 7683.428400000001
$ws.Range("I132").Value = 2561.1428
$ws.Range("N135").Value = -160189.5
$ws.Range("J135").Value = 150049.5
$ws.Range("L135").Value = 150049.5
$ws.Range("H135").Value = 150049.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L22").Value = 500000000
$ws.Range("N22").Value = -500000590
$ws.Range("J22").Value = 500000000
$ws.Range("H22").Value = 500000000
$ws.Range("K22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("L27").Value = 500000000
$ws.Range("N27").Value = -500000214
$ws.Range("H27").Value = 500000000
$ws.Range("J27").Value = 500000000
$ws.Range("K27").Value = 0
$ws.Range("J55").Value = 1365.6
$ws.Range("H55").Value = 1145.32
$ws.Range("L55").Value = 1365.6
$ws.Range("N55").Value = -1711.6
$ws.Range("J68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("N82").Value = -2184.1111
$ws.Range("M82").Value = -549.8
$ws.Range("L82").Value = 1462.1111
$ws.Range("J82").Value = 1462.1111
$ws.Range("H82").Value = 1171.9474
$ws.Range("K82").Value = 910.8
$ws.Range("I82").Value = 910.8
$ws.Range("K85").Value = 910.8
$ws.Range("I85").Value = 910.8
$ws.Range("N85").Value = -3958.1111
$ws.Range("M85").Value = 337.2
$ws.Range("L85").Value = 1462.1111
$ws.Range("J85").Value = 1462.1111
$ws.Range("H85").Value = 1171.9474
$ws.Range("L98").Value = 60461.5
$ws.Range("N98").Value = -66451.5
$ws.Range("J98").Value = 60461.5
$ws.Range("H98").Value = 60461.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I5").Value = 2000
$ws.Range("H5").Value = 2000
$ws.Range("M5").Value = -1888
$ws.Range("K5").Value = 2000
$ws.Range("H49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("H62").Value = 9995
$ws.Range("H65").Value = 9995
$ws.Range("K65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J81").Value = 12957.889
$ws.Range("H81").Value = 6472.069
$ws.Range("M81").Value = -6045.9
$ws.Range("K81").Value = 7106.9
$ws.Range("I81").Value = 3553.45
$ws.Range("L81").Value = 25915.778
$ws.Range("N81").Value = -28037.778
$ws.Range("J84").Value = 12957.889
$ws.Range("H84").Value = 6472.069
$ws.Range("K84").Value = 35534.5
$ws.Range("I84").Value = 3553.45
$ws.Range("M84").Value = -30230.5
$ws.Range("L84").Value = 129578.89
$ws.Range("N84").Value = -140186.89
$ws.Range("K107").Value = 1005
$ws.Range("I107").Value = 335
$ws.Range("M107").Value = 915
$ws.Range("H107").Value = 335
$ws.Range("M132").Value = -2864.428400000001
$ws.Range("H132").Value = 90910824
$ws.Range("K132").Value = 5394.428400000001
$ws.Range("I132").Value = 1798.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M22").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("M27").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M49").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
